$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 495.7143
$ws.Cells.Item(41, 9).Value = 98
$ws.Cells.Item(41, 11).Value = 98
$ws.Cells.Item(41, 13).Value = 342
$ws.Cells.Item(53, 8).Value = 146.59259
$ws.Cells.Item(53, 9).Value = 98.61539
$ws.Cells.Item(53, 10).Value = 191.14285
$ws.Cells.Item(53, 11).Value = 98.61539
$ws.Cells.Item(53, 12).Value = 191.14285
$ws.Cells.Item(53, 13).Value = 538.38461
$ws.Cells.Item(53, 14).Value = -1465.14285
$ws.Cells.Item(106, 8).Value = 1751.2222
$ws.Cells.Item(106, 9).Value = 935
$ws.Cells.Item(106, 10).Value = 2159.3333
$ws.Cells.Item(106, 11).Value = 935
$ws.Cells.Item(106, 12).Value = 2159.3333
$ws.Cells.Item(106, 13).Value = -304
$ws.Cells.Item(106, 14).Value = -3421.3333
$ws.Cells.Item(133, 8).Value = 53800
$ws.Cells.Item(133, 10).Value = 53800
$ws.Cells.Item(133, 12).Value = 53800
$ws.Cells.Item(133, 14).Value = -63920
$ws.Cells.Item(137, 8).Value = 1542.877
$ws.Cells.Item(137, 9).Value = 1156.8223
$ws.Cells.Item(137, 10).Value = 2411.5
$ws.Cells.Item(137, 11).Value = 3470.4669
$ws.Cells.Item(137, 12).Value = 7234.5
$ws.Cells.Item(137, 13).Value = -920.4669000000004
$ws.Cells.Item(137, 14).Value = -12334.5
$ws.Cells.Item(141, 8).Value = 2445.5454
$ws.Cells.Item(141, 9).Value = 1342.561
$ws.Cells.Item(141, 10).Value = 5675.7144
$ws.Cells.Item(141, 11).Value = 4027.683
$ws.Cells.Item(141, 12).Value = 17027.1432
$ws.Cells.Item(141, 13).Value = 1152.317
$ws.Cells.Item(141, 14).Value = -27387.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 6840.827
$ws.Cells.Item(61, 9).Value = 3704.4546
$ws.Cells.Item(61, 10).Value = 24090.875
$ws.Cells.Item(61, 11).Value = 3704.4546
$ws.Cells.Item(61, 12).Value = 24090.875
$ws.Cells.Item(61, 13).Value = -3492.4546
$ws.Cells.Item(61, 14).Value = -24514.875
$ws.Cells.Item(112, 8).Value = 39387
$ws.Cells.Item(112, 10).Value = 39387
$ws.Cells.Item(112, 12).Value = 39387
$ws.Cells.Item(112, 14).Value = -42341
$ws.Cells.Item(122, 8).Value = 1687.6086
$ws.Cells.Item(122, 9).Value = 1302.5333
$ws.Cells.Item(122, 10).Value = 2409.625
$ws.Cells.Item(122, 11).Value = 3907.5999
$ws.Cells.Item(122, 12).Value = 7228.875
$ws.Cells.Item(122, 13).Value = -1457.5999
$ws.Cells.Item(122, 14).Value = -12128.875
$ws.Cells.Item(123, 8).Value = 45933.332
$ws.Cells.Item(123, 10).Value = 45933.332
$ws.Cells.Item(123, 12).Value = 45933.332
$ws.Cells.Item(123, 14).Value = -55733.332
$ws.Cells.Item(136, 8).Value = 6840.827
$ws.Cells.Item(136, 9).Value = 3704.4546
$ws.Cells.Item(136, 10).Value = 24090.875
$ws.Cells.Item(136, 11).Value = 11113.3638
$ws.Cells.Item(136, 12).Value = 72272.625
$ws.Cells.Item(136, 13).Value = -8563.363799999999
$ws.Cells.Item(136, 14).Value = -77372.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 311.2143
$ws.Cells.Item(22, 9).Value = 296.69232
$ws.Cells.Item(22, 11).Value = 296.69232
$ws.Cells.Item(22, 13).Value = -123.69232
$ws.Cells.Item(86, 8).Value = 1735.7142
$ws.Cells.Item(86, 9).Value = 1741.6666
$ws.Cells.Item(86, 10).Value = 1700
$ws.Cells.Item(86, 11).Value = 1741.6666
$ws.Cells.Item(86, 12).Value = 1700
$ws.Cells.Item(86, 13).Value = -618.6666
$ws.Cells.Item(86, 14).Value = -3946
$ws.Cells.Item(89, 8).Value = 1735.7142
$ws.Cells.Item(89, 9).Value = 1741.6666
$ws.Cells.Item(89, 10).Value = 1700
$ws.Cells.Item(89, 11).Value = 8708.333000000001
$ws.Cells.Item(89, 12).Value = 8500
$ws.Cells.Item(89, 13).Value = -3092.333000000001
$ws.Cells.Item(89, 14).Value = -19732
$ws.Cells.Item(107, 8).Value = 1288.8889
$ws.Cells.Item(107, 9).Value = 1100
$ws.Cells.Item(107, 11).Value = 1100
$ws.Cells.Item(107, 13).Value = 820
$ws.Cells.Item(132, 8).Value = 49861.54
$ws.Cells.Item(132, 10).Value = 58820
$ws.Cells.Item(132, 12).Value = 58820
$ws.Cells.Item(132, 14).Value = -68940

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 139.13043
$ws.Cells.Item(7, 9).Value = 121.53846
$ws.Cells.Item(7, 10).Value = 162
$ws.Cells.Item(7, 11).Value = 121.53846
$ws.Cells.Item(7, 12).Value = 162
$ws.Cells.Item(7, 13).Value = -8.538460000000001
$ws.Cells.Item(7, 14).Value = -388
$ws.Cells.Item(31, 8).Value = 2077.3215
$ws.Cells.Item(31, 9).Value = 1460.7073
$ws.Cells.Item(31, 10).Value = 3762.7334
$ws.Cells.Item(31, 11).Value = 1460.7073
$ws.Cells.Item(31, 12).Value = 3762.7334
$ws.Cells.Item(31, 13).Value = -1165.7073
$ws.Cells.Item(31, 14).Value = -4352.7334
$ws.Cells.Item(34, 8).Value = 2077.3215
$ws.Cells.Item(34, 9).Value = 1460.7073
$ws.Cells.Item(34, 10).Value = 3762.7334
$ws.Cells.Item(34, 11).Value = 1460.7073
$ws.Cells.Item(34, 12).Value = 3762.7334
$ws.Cells.Item(34, 13).Value = -1258.7073
$ws.Cells.Item(34, 14).Value = -4166.7334
$ws.Cells.Item(99, 8).Value = 5189.125
$ws.Cells.Item(99, 9).Value = 4649.5
$ws.Cells.Item(99, 10).Value = 5369
$ws.Cells.Item(99, 11).Value = 4649.5
$ws.Cells.Item(99, 12).Value = 5369
$ws.Cells.Item(99, 13).Value = -3151.5
$ws.Cells.Item(99, 14).Value = -8365
$ws.Cells.Item(105, 8).Value = 1383.6428
$ws.Cells.Item(105, 9).Value = 825
$ws.Cells.Item(105, 11).Value = 825
$ws.Cells.Item(105, 13).Value = 922
$ws.Cells.Item(122, 8).Value = 10571.429
$ws.Cells.Item(122, 9).Value = 12000
$ws.Cells.Item(122, 10).Value = 7000
$ws.Cells.Item(122, 11).Value = 36000
$ws.Cells.Item(122, 12).Value = 21000
$ws.Cells.Item(122, 13).Value = -33550
$ws.Cells.Item(122, 14).Value = -25900
$ws.Cells.Item(126, 8).Value = 5189.125
$ws.Cells.Item(126, 9).Value = 4649.5
$ws.Cells.Item(126, 10).Value = 5369
$ws.Cells.Item(126, 11).Value = 13948.5
$ws.Cells.Item(126, 12).Value = 16107
$ws.Cells.Item(126, 13).Value = -11478.5
$ws.Cells.Item(126, 14).Value = -21047
$ws.Cells.Item(132, 8).Value = 2619.8215
$ws.Cells.Item(132, 9).Value = 2821.6558
$ws.Cells.Item(132, 10).Value = 2084.5217
$ws.Cells.Item(132, 11).Value = 8464.9674
$ws.Cells.Item(132, 12).Value = 6253.5651
$ws.Cells.Item(132, 13).Value = -5934.9674
$ws.Cells.Item(132, 14).Value = -11313.5651

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(8, 8).Value = 62.285713
$ws.Cells.Item(8, 9).Value = 62.285713
$ws.Cells.Item(8, 11).Value = 186.857139
$ws.Cells.Item(8, 13).Value = -47.85713900000002
$ws.Cells.Item(12, 8).Value = 29411976
$ws.Cells.Item(12, 9).Value = 62500176
$ws.Cells.Item(12, 10).Value = 243.11111
$ws.Cells.Item(12, 11).Value = 187500528
$ws.Cells.Item(12, 12).Value = 729.3333299999999
$ws.Cells.Item(12, 13).Value = -187500355
$ws.Cells.Item(12, 14).Value = -1075.33333
$ws.Cells.Item(14, 8).Value = 37256.32
$ws.Cells.Item(14, 9).Value = 37256.32
$ws.Cells.Item(14, 11).Value = 111768.96
$ws.Cells.Item(14, 13).Value = -111595.96
$ws.Cells.Item(92, 8).Value = 473.0909
$ws.Cells.Item(92, 10).Value = 487
$ws.Cells.Item(92, 12).Value = 1461
$ws.Cells.Item(92, 14).Value = -3957

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 6348.375
$ws.Cells.Item(70, 9).Value = 6014.8
$ws.Cells.Item(70, 10).Value = 6500
$ws.Cells.Item(70, 11).Value = 6014.8
$ws.Cells.Item(70, 12).Value = 6500
$ws.Cells.Item(70, 13).Value = -5744.8
$ws.Cells.Item(70, 14).Value = -7040
$ws.Cells.Item(73, 8).Value = 6348.375
$ws.Cells.Item(73, 9).Value = 6014.8
$ws.Cells.Item(73, 10).Value = 6500
$ws.Cells.Item(73, 11).Value = 6014.8
$ws.Cells.Item(73, 12).Value = 6500
$ws.Cells.Item(73, 13).Value = -5078.8
$ws.Cells.Item(73, 14).Value = -8372
$ws.Cells.Item(97, 8).Value = 1523.3334
$ws.Cells.Item(97, 9).Value = 1458.5714
$ws.Cells.Item(97, 10).Value = 1750
$ws.Cells.Item(97, 11).Value = 1458.5714
$ws.Cells.Item(97, 12).Value = 1750
$ws.Cells.Item(97, 13).Value = -962.5714
$ws.Cells.Item(97, 14).Value = -2742
$ws.Cells.Item(111, 8).Value = 25599.6
$ws.Cells.Item(111, 10).Value = 25599.6
$ws.Cells.Item(111, 12).Value = 25599.6
$ws.Cells.Item(111, 14).Value = -31733.6
$ws.Cells.Item(132, 8).Value = 5925.659
$ws.Cells.Item(132, 9).Value = 4071.4473
$ws.Cells.Item(132, 10).Value = 17669
$ws.Cells.Item(132, 11).Value = 12214.3419
$ws.Cells.Item(132, 12).Value = 53007
$ws.Cells.Item(132, 13).Value = -9684.341899999999
$ws.Cells.Item(132, 14).Value = -58067

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 973.6
$ws.Cells.Item(22, 9).Value = 1048.4
$ws.Cells.Item(22, 10).Value = 936.2
$ws.Cells.Item(22, 11).Value = 1048.4
$ws.Cells.Item(22, 12).Value = 936.2
$ws.Cells.Item(22, 13).Value = -753.4000000000001
$ws.Cells.Item(22, 14).Value = -1526.2
$ws.Cells.Item(27, 8).Value = 973.6
$ws.Cells.Item(27, 9).Value = 1048.4
$ws.Cells.Item(27, 10).Value = 936.2
$ws.Cells.Item(27, 11).Value = 1048.4
$ws.Cells.Item(27, 12).Value = 936.2
$ws.Cells.Item(27, 13).Value = -941.4000000000001
$ws.Cells.Item(27, 14).Value = -1150.2
$ws.Cells.Item(46, 8).Value = 3000
$ws.Cells.Item(46, 9).Value = 3000
$ws.Cells.Item(46, 10).Value = 0
$ws.Cells.Item(46, 11).Value = 3000
$ws.Cells.Item(46, 12).Value = 0
$ws.Cells.Item(46, 13).Value = -2812
$ws.Cells.Item(46, 14).ClearContents()
$ws.Cells.Item(61, 8).Value = 857447.7
$ws.Cells.Item(61, 9).Value = 23936.7
$ws.Cells.Item(61, 11).Value = 23936.7
$ws.Cells.Item(61, 13).Value = -23734.7
$ws.Cells.Item(68, 8).Value = 900
$ws.Cells.Item(68, 9).Value = 800
$ws.Cells.Item(68, 10).Value = 1000
$ws.Cells.Item(68, 11).Value = 800
$ws.Cells.Item(68, 12).Value = 1000
$ws.Cells.Item(68, 13).Value = -51
$ws.Cells.Item(68, 14).Value = -2498
$ws.Cells.Item(71, 8).Value = 900
$ws.Cells.Item(71, 9).Value = 800
$ws.Cells.Item(71, 10).Value = 1000
$ws.Cells.Item(71, 11).Value = 4000
$ws.Cells.Item(71, 12).Value = 5000
$ws.Cells.Item(71, 13).Value = -256
$ws.Cells.Item(71, 14).Value = -12488
$ws.Cells.Item(93, 8).Value = 1655.3
$ws.Cells.Item(93, 9).Value = 1450.3334
$ws.Cells.Item(93, 10).Value = 3500
$ws.Cells.Item(93, 11).Value = 1450.3334
$ws.Cells.Item(93, 12).Value = 3500
$ws.Cells.Item(93, 13).Value = -202.3334
$ws.Cells.Item(93, 14).Value = -5996
$ws.Cells.Item(113, 8).Value = 857447.7
$ws.Cells.Item(113, 9).Value = 23936.7
$ws.Cells.Item(113, 11).Value = 23936.7
$ws.Cells.Item(113, 13).Value = -21766.7
$ws.Cells.Item(122, 8).Value = 6405.523
$ws.Cells.Item(122, 9).Value = 6008.841
$ws.Cells.Item(122, 11).Value = 18026.523
$ws.Cells.Item(122, 13).Value = -15576.523

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 5400
$ws.Cells.Item(14, 10).Value = 5400
$ws.Cells.Item(14, 12).Value = 5400
$ws.Cells.Item(14, 14).Value = -5736
$ws.Cells.Item(22, 8).Value = 13
$ws.Cells.Item(22, 9).Value = 13
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 13
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 280
$ws.Cells.Item(22, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 2667.75
$ws.Cells.Item(122, 9).Value = 1693.0714
$ws.Cells.Item(122, 10).Value = 3642.4285
$ws.Cells.Item(122, 11).Value = 5079.2142
$ws.Cells.Item(122, 12).Value = 10927.2855
$ws.Cells.Item(122, 13).Value = -2629.2142
$ws.Cells.Item(122, 14).Value = -15827.2855
$ws.Cells.Item(132, 8).Value = 2129.84
$ws.Cells.Item(132, 9).Value = 755.4
$ws.Cells.Item(132, 10).Value = 3046.1333
$ws.Cells.Item(132, 11).Value = 2266.2
$ws.Cells.Item(132, 12).Value = 9138.3999
$ws.Cells.Item(132, 13).Value = 263.8000000000002
$ws.Cells.Item(132, 14).Value = -14198.3999
$ws.Cells.Item(136, 8).Value = 5071.4614
$ws.Cells.Item(136, 9).Value = 3475.6047
$ws.Cells.Item(136, 10).Value = 8190.636
$ws.Cells.Item(136, 11).Value = 10426.8141
$ws.Cells.Item(136, 12).Value = 24571.908
$ws.Cells.Item(136, 13).Value = -7876.8141
$ws.Cells.Item(136, 14).Value = -29671.908
